$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 704, pushing existing rows 704:759 down to 706:761
$ws.Rows("704:705").Insert()

# Populate the two newly inserted rows (704 and 705) with new data
# Row 704
$ws.Range("A704").Value = 7
$ws.Range("B704").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C704").Value = "Ñuble"
$ws.Range("D704").Value = 45223
$ws.Range("E704").Value = 16
$ws.Range("F704").Value = 100114001
$ws.Range("G704").Value = "Papa"
$ws.Range("H704").Value = "Asterix"
$ws.Range("I704").Value = "1a (guarda)"
$ws.Range("J704").Value = 270
$ws.Range("K704").Value = 28000
$ws.Range("L704").Value = 30000
$ws.Range("M704").Value = 28889
$ws.Range("N704").Value = "$/saco 25 kilos"
$ws.Range("O704").Value = "Región de Los Lagos"
$ws.Range("P704").Value = 1156
$ws.Range("Q704").Value = 25
$ws.Range("R704").Value = "Hortaliza"

# Row 705
$ws.Range("A705").Value = 7
$ws.Range("B705").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C705").Value = "Ñuble"
$ws.Range("D705").Value = 45223
$ws.Range("E705").Value = 16
$ws.Range("F705").Value = 100114001
$ws.Range("G705").Value = "Papa"
$ws.Range("H705").Value = "Rodeo"
$ws.Range("I705").Value = "1a (guarda)"
$ws.Range("J705").Value = 150
$ws.Range("K705").Value = 28000
$ws.Range("L705").Value = 28000
$ws.Range("M705").Value = 28000
$ws.Range("N705").Value = "$/saco 25 kilos"
$ws.Range("O705").Value = "Región de Los Lagos"
$ws.Range("P705").Value = 1120
$ws.Range("Q705").Value = 25
$ws.Range("R705").Value = "Hortaliza"
